# Change the subtitle text on the two "Pitfalls to Avoid" slides.
#
# Slide 1 (title slide): subtitle paragraph "Pitfalls to Avoid" becomes two
# runs: "And " followed by "avoiding pitfalls while doing it".
#
# Slide 2 (agenda/section slide): subtitle run "Pitfalls to Avoid" (bold)
# becomes "And avoiding pitfalls while doing it" (no longer bold).

$p = $ppt.ActivePresentation

# --- Slide 1 -----------------------------------------------------------
$s1 = $p.Slides.Item(1)
$sh1 = $s1.Shapes.Item("Subtitle 4")   # subtitle placeholder
$tr1 = $sh1.TextFrame.TextRange
$para1 = $tr1.Paragraphs(2)            # paragraph holding "Pitfalls to Avoid"
$run1 = $para1.Runs(1)

# Insert a brand-new leading run ("And ") ahead of the existing run, then
# rewrite the existing run's text in place so its own formatting (and the
# dirty="0" marker already on it) is preserved.
[void]$run1.InsertBefore("And ")
$para1.Runs(2).Text = "avoiding pitfalls while doing it"

# --- Slide 2 -----------------------------------------------------------
$s2 = $p.Slides.Item(2)
$sh2 = $s2.Shapes.Item("Subtitle 2")   # subtitle placeholder
$tr2 = $sh2.TextFrame.TextRange
$para2 = $tr2.Paragraphs(1)            # paragraph holding "Pitfalls to Avoid"
$run2 = $para2.Runs(1)

$run2.Text = "And avoiding pitfalls while doing it"
$run2.Font.Bold = $false
